# Update the cryptocurrency price/volume snapshot (GitHub Actions refresh).
# Price/percent cells are stored as plain text in the workbook, so values
# that look like plain numbers ("249.60", "0.0808", ...) are written with a
# leading apostrophe to force Excel to keep them as text (preserving exact
# formatting / trailing zeros instead of silently coercing to a float), and
# the cell style is then reset to "Normal" so no stray numeric formatting
# is left behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.018.23'
$ws.Range("E2").Value = '  +1.08%  '
$ws.Range("D3").Value = '2.062.50'
$ws.Range("E3").Value = '  -1.44%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").Value = "'249.60"
$ws.Range("E5").Value = '  -1.24%  '
$ws.Range("D6").Value = "'0.673"
$ws.Range("E6").Value = '  +1.95%  '
$ws.Range("E7").Value = '  -0.01%  '
$ws.Range("D8").Value = "'55.11"
$ws.Range("E8").Value = '  +11.52%  '
$ws.Range("D9").Value = "'60.72"
$ws.Range("E9").Value = '  +0.37%  '
$ws.Range("D10").Value = "'0.383"
$ws.Range("E10").Value = '  +1.95%  '
$ws.Range("D11").Value = "'0.0808"
$ws.Range("E11").Value = '  +8.09%  '
$ws.Range("E12").Value = '  +5.87%  '
$ws.Range("D13").Value = "'15.09"
$ws.Range("E13").Value = '  +1.55%  '
$ws.Range("D14").Value = '2.364.71'
$ws.Range("E14").Value = '  -1.42%  '
$ws.Range("D15").Value = "'0.819"
$ws.Range("E15").Value = '  -2.44%  '
$ws.Range("E16").Value = '  +3.44%  '
$ws.Range("D17").Value = '2.066.58'
$ws.Range("E17").Value = '  -1.32%  '
$ws.Range("D18").Value = '36.990.76'
$ws.Range("E18").Value = '  +1.12%  '
$ws.Range("E19").Value = '  +13.47%  '
$ws.Range("D20").Value = "'73.51"
$ws.Range("E20").Value = '  +0.37%  '
$ws.Range("D21").Value = "'14.24"
$ws.Range("E21").Value = '  +7.08%  '
$ws.Range("D22").Value = "'5.40"
$ws.Range("E22").Value = '  +1.36%  '
$ws.Range("D23").Value = "'237.71"
$ws.Range("E23").Value = '  -1.17%  '
$ws.Range("E24").Value = '  -0.11%  '
$ws.Range("E25").Value = '  -4.36%  '
$ws.Range("D26").Value = "'174.27"
$ws.Range("E26").Value = '  +1.72%  '
$ws.Range("D27").Value = "'9.14"
$ws.Range("E27").Value = '  -1.65%  '
$ws.Range("E28").Value = '  -4.92%  '
$ws.Range("E29").Value = '  +0.13%  '
$ws.Range("E30").Value = '  +1.99%  '
$ws.Range("D31").Value = "'4.61"
$ws.Range("E31").Value = '  +2.16%  '
$ws.Range("E32").Value = '  +6.52%  '
$ws.Range("D33").Value = "'0.0629"
$ws.Range("E33").Value = '  +2.01%  '
$ws.Range("E34").Value = '  +7.27%  '
$ws.Range("D35").Value = "'0.0894"
$ws.Range("E35").Value = '  -1.79%  '
$ws.Range("E36").Value = '  +0.05%  '
$ws.Range("E37").Value = '  -6.18%  '
$ws.Range("E38").Value = '  -5.01%  '
$ws.Range("E39").Value = '  -0.05%  '
$ws.Range("E40").Value = '  +24.67%  '
$ws.Range("D41").Value = "'17.98"
$ws.Range("E41").Value = '  +7.96%  '
$ws.Range("E42").Value = '  +0.42%  '
$ws.Range("E43").Value = '  -1.94%  '
$ws.Range("D44").Value = "'96.85"
$ws.Range("E44").Value = '  -1.31%  '
$ws.Range("E45").Value = '  +0.46%  '
$ws.Range("B46").Value = 'Gas'
$ws.Range("C46").Value = 'https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas'
$ws.Range("D46").Value = "'14.09"
$ws.Range("E46").Value = '  -50.34%  '
$ws.Range("B47").Value = 'FTXToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D47").Value = "'4.11"
$ws.Range("E47").Value = '  +35.59%  '
$ws.Range("D48").Value = "'2.42"
$ws.Range("E48").Value = '  +6.70%  '
$ws.Range("D49").Value = "'4.28"
$ws.Range("E49").Value = '  +10.53%  '
$ws.Range("D50").Value = '1.302.32'
$ws.Range("E50").Value = '  -2.66%  '
$ws.Range("E51").Value = '  +1.40%  '

# Clear the forced-text "quote prefix" styling introduced above so the
# affected cells fall back to the workbook's default (unstyled) format.
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
